# Update "paises" (Covid-19 countries) worksheet:
#   - refresh the "Datos actualizados..." timestamp in A1
#   - refresh case counters for several countries (new data snapshot)
#   - the table is kept sorted descending by "Casos totales" (col B),
#     so a few rows swap the country name in col A as their neighbours
#     re-rank following the data refresh (e.g. Guinea overtakes
#     Uzbekistan & Senegal, Costa Rica overtakes Sudan del Sur, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / last-updated timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 28 de Mayo de 2020 a las 22:10"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1760816
$ws.Cells.Item(4, 3).Value = 15013
$ws.Cells.Item(4, 5).Value = 1163103
$ws.Cells.Item(4, 7).Value = 882
$ws.Cells.Item(4, 8).Value = 102989

# Row 5: Brasil
$ws.Cells.Item(5, 2).Value = 419340
$ws.Cells.Item(5, 3).Value = 4679
$ws.Cells.Item(5, 4).Value = 192302
$ws.Cells.Item(5, 5).Value = 201093
$ws.Cells.Item(5, 7).Value = 248
$ws.Cells.Item(5, 8).Value = 25945

# Row 11: Alemania
$ws.Cells.Item(11, 2).Value = 182452
$ws.Cells.Item(11, 3).Value = 557
$ws.Cells.Item(11, 5).Value = 10682
$ws.Cells.Item(11, 7).Value = 37
$ws.Cells.Item(11, 8).Value = 8570

# Row 12: India
$ws.Cells.Item(12, 2).Value = 165386
$ws.Cells.Item(12, 3).Value = 7300
$ws.Cells.Item(12, 4).Value = 70920
$ws.Cells.Item(12, 5).Value = 89755

# Row 16: Canada
$ws.Cells.Item(16, 2).Value = 88473
$ws.Cells.Item(16, 3).Value = 954
$ws.Cells.Item(16, 4).Value = 46768
$ws.Cells.Item(16, 5).Value = 34832

# Row 42: Israel
$ws.Cells.Item(42, 2).Value = 16872
$ws.Cells.Item(42, 3).Value = 79
$ws.Cells.Item(42, 4).Value = 14679
$ws.Cells.Item(42, 5).Value = 1909
$ws.Cells.Item(42, 7).Value = 3
$ws.Cells.Item(42, 8).Value = 284

# Row 77: Guinea
$ws.Cells.Item(77, 1).Value = "Guinea"
$ws.Cells.Item(77, 2).Value = 3553
$ws.Cells.Item(77, 3).Value = 278
$ws.Cells.Item(77, 4).Value = 1950
$ws.Cells.Item(77, 5).Value = 1581
$ws.Cells.Item(77, 7).Value = 2
$ws.Cells.Item(77, 8).Value = 22

# Row 78: Uzbekistan
$ws.Cells.Item(78, 1).Value = "Uzbekistan"
$ws.Cells.Item(78, 2).Value = 3444
$ws.Cells.Item(78, 3).Value = 75
$ws.Cells.Item(78, 4).Value = 2694
$ws.Cells.Item(78, 5).Value = 736
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 14

# Row 79: Senegal
$ws.Cells.Item(79, 1).Value = "Senegal"
$ws.Cells.Item(79, 2).Value = 3348
$ws.Cells.Item(79, 3).Value = 95
$ws.Cells.Item(79, 4).Value = 1686
$ws.Cells.Item(79, 5).Value = 1623
$ws.Cells.Item(79, 7).Value = 1
$ws.Cells.Item(79, 8).Value = 39

# Row 84: Costa de Marfil
$ws.Cells.Item(84, 2).Value = 2641
$ws.Cells.Item(84, 3).Value = 85
$ws.Cells.Item(84, 4).Value = 1326
$ws.Cells.Item(84, 5).Value = 1283
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = 32

# Row 116: Costa Rica
$ws.Cells.Item(116, 1).Value = "Costa Rica"
$ws.Cells.Item(116, 2).Value = 1000
$ws.Cells.Item(116, 3).Value = 16
$ws.Cells.Item(116, 4).Value = 646
$ws.Cells.Item(116, 5).Value = 344

# Row 117: Sudan del Sur
$ws.Cells.Item(117, 1).Value = "Sudan del Sur"
$ws.Cells.Item(117, 2).Value = 994
$ws.Cells.Item(117, 4).Value = 6
$ws.Cells.Item(117, 5).Value = 978

# Row 119: Republica de Chipre
$ws.Cells.Item(119, 4).Value = 784
$ws.Cells.Item(119, 5).Value = 140

# Row 140: Estado de Palestina
$ws.Cells.Item(140, 1).Value = "Estado de Palestina"
$ws.Cells.Item(140, 2).Value = 446
$ws.Cells.Item(140, 3).Value = 12
$ws.Cells.Item(140, 4).Value = 365
$ws.Cells.Item(140, 5).Value = 78
$ws.Cells.Item(140, 8).Value = 3

# Row 141: Taiwan
$ws.Cells.Item(141, 1).Value = "Taiwan"
$ws.Cells.Item(141, 2).Value = 441
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = 420
$ws.Cells.Item(141, 5).Value = 14
$ws.Cells.Item(141, 8).Value = 7

# Row 142: Togo
$ws.Cells.Item(142, 1).Value = "Togo"
$ws.Cells.Item(142, 2).Value = 422
$ws.Cells.Item(142, 3).Value = 27
$ws.Cells.Item(142, 4).Value = 197
$ws.Cells.Item(142, 5).Value = 212
$ws.Cells.Item(142, 8).Value = 13

# Row 143: Guayana Francesa
$ws.Cells.Item(143, 1).Value = "Guayana Francesa"
$ws.Cells.Item(143, 2).Value = 406
$ws.Cells.Item(143, 4).Value = 150
$ws.Cells.Item(143, 5).Value = 255
$ws.Cells.Item(143, 8).Value = 1

# Row 198: Curazao
$ws.Cells.Item(198, 1).Value = "Curazao"
$ws.Cells.Item(198, 4).Value = 14
$ws.Cells.Item(198, 8).Value = 1

# Row 199: Fiyi
$ws.Cells.Item(199, 1).Value = "Fiyi"
$ws.Cells.Item(199, 4).Value = 15
$ws.Cells.Item(199, 8).Value = 0

# Row 200: Santa Lucia
$ws.Cells.Item(200, 1).Value = "Santa Lucia"
$ws.Cells.Item(200, 4).Value = 18
$ws.Cells.Item(200, 8).Value = 0

# Row 201: Belice
$ws.Cells.Item(201, 1).Value = "Belice"
$ws.Cells.Item(201, 4).Value = 16
$ws.Cells.Item(201, 8).Value = 2

# Row 213: Papua Nueva Guinea
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 8).Value = 0

# Row 214: Islas Virgenes Britanicas
$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 8).Value = 1
